$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) & data row (row 2) ----
# Values are written in the same order the original author typed them in
# (row 1 left-to-right first, then a couple of row-2 cells, then A1,
# then the rest of row 2) so the shared-string table comes out in the
# same sequence.
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Branch Name"
$ws.Range("D1").Value = "Year"
$ws.Range("E1").Value = "Section"
$ws.Range("F1").Value = "Discipline ID"
$ws.Range("G1").Value = "Semester"
$ws.Range("H1").Value = "Subject 1"
$ws.Range("I1").Value = "Subject 2"
$ws.Range("J1").Value = "Subject 3"
$ws.Range("K1").Value = "Subject 4"
$ws.Range("L1").Value = "Subject 5"
$ws.Range("C2").Value = "CSE"
$ws.Range("A1").Value = "unique_id"
$ws.Range("B2").Value = "dummy1"

$ws.Range("A2").Value = 21111111
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 24
$ws.Range("I2").Value = 26

# ---- Formatting ----
# Clear the stray bold-ish style that used to live on A1 so it goes back
# to the default (unstyled) cell format.
$ws.Range("A1").ClearFormats()

# Header row: centered, middle-aligned, wrapped text.
$headerRng = $ws.Range("B1:L1")
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4108
$headerRng.WrapText = $true

# Row height for the (now two-line) header row.
$ws.Rows.Item(1).RowHeight = 28.8

# Column widths.
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 24.21875

# Selection / active cell.
$ws.Range("A2").Select()

# Theme accent colors: swap accent1 <-> accent5.
$scheme = $wb.Theme.ThemeColorScheme
$scheme.Colors(5).RGB = 13998939   # 5B9BD5
$scheme.Colors(9).RGB = 12874308   # 4472C4
